$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of updated Price values are unpadded decimals that Excel's
# COM type-inference would otherwise coerce to numbers (losing trailing
# zeros, e.g. "1.00" -> 1). Mark just those cells as Text first so the
# literal digit string round-trips exactly, matching the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.652.03'
$ws.Cells.Item(2, 5).Value = '  -0.09%  '
$ws.Cells.Item(3, 4).Value = '1.530.09'
$ws.Cells.Item(3, 5).Value = '  -1.55%  '
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$ws.Cells.Item(5, 4).Value = '205.68'
$ws.Cells.Item(5, 5).Value = '  +0.10%  '
$ws.Cells.Item(6, 4).Value = '0.483'
$ws.Cells.Item(6, 5).Value = '  -1.10%  '
$ws.Cells.Item(7, 5).Value = '  -0.07%  '
$ws.Cells.Item(8, 4).Value = '21.27'
$ws.Cells.Item(8, 5).Value = '  -2.23%  '
$ws.Cells.Item(9, 5).Value = '  -1.17%  '
$ws.Cells.Item(10, 5).Value = '  -0.56%  '
$ws.Cells.Item(11, 5).Value = '  -1.23%  '
$ws.Cells.Item(12, 4).Value = '1.748.57'
$ws.Cells.Item(12, 5).Value = '  -1.58%  '
$ws.Cells.Item(13, 4).Value = '1.527.15'
$ws.Cells.Item(13, 5).Value = '  -1.86%  '
$ws.Cells.Item(14, 4).Value = '3.67'
$ws.Cells.Item(14, 5).Value = '  -1.33%  '
$ws.Cells.Item(15, 4).Value = '0.504'
$ws.Cells.Item(15, 5).Value = '  -1.09%  '
$ws.Cells.Item(16, 5).Value = '  -0.13%  '
$ws.Cells.Item(17, 4).Value = '26.657.36'
$ws.Cells.Item(17, 5).Value = '  -0.27%  '
$ws.Cells.Item(18, 4).Value = '212.41'
$ws.Cells.Item(18, 5).Value = '  -0.47%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0683'
$ws.Cells.Item(19, 5).Value = '  +1.31%  '
$ws.Cells.Item(20, 5).Value = '  -1.74%  '
$ws.Cells.Item(21, 5).Value = '  -0.10%  '
$ws.Cells.Item(22, 5).Value = '  -1.72%  '
$ws.Cells.Item(23, 4).Value = '9.08'
$ws.Cells.Item(23, 5).Value = '  -2.92%  '
$ws.Cells.Item(24, 5).Value = '  -3.85%  '
$ws.Cells.Item(25, 4).Value = '152.27'
$ws.Cells.Item(25, 5).Value = '  -0.31%  '
$ws.Cells.Item(26, 4).Value = '6.52'
$ws.Cells.Item(26, 5).Value = '  -3.09%  '
$ws.Cells.Item(27, 4).Value = '14.81'
$ws.Cells.Item(27, 5).Value = '  +0.00%  '
$ws.Cells.Item(28, 5).Value = '  -0.15%  '
$ws.Cells.Item(29, 5).Value = '  -0.71%  '
$ws.Cells.Item(30, 4).Value = '1.09'
$ws.Cells.Item(30, 5).Value = '  -1.08%  '
$ws.Cells.Item(31, 5).Value = '  -1.88%  '
$ws.Cells.Item(32, 5).Value = '  +2.83%  '
$ws.Cells.Item(33, 4).Value = '1.353.45'
$ws.Cells.Item(33, 5).Value = '  -1.76%  '
$ws.Cells.Item(34, 5).Value = '  +0.32%  '
$ws.Cells.Item(35, 5).Value = '  -3.36%  '
$ws.Cells.Item(36, 4).Value = '0.951'
$ws.Cells.Item(36, 5).Value = '  +1.52%  '
$ws.Cells.Item(37, 5).Value = '  -0.68%  '
$ws.Cells.Item(38, 5).Value = '  +0.36%  '
$ws.Cells.Item(39, 4).Value = '0.521'
$ws.Cells.Item(39, 5).Value = '  +0.76%  '
$ws.Cells.Item(40, 2).Value = 'PaxDollar'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  -0.11%  '
$ws.Cells.Item(41, 2).Value = 'ARBITRUM'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(41, 4).Value = '0.795'
$ws.Cells.Item(41, 5).Value = '  -1.45%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(42, 4).Value = '5.68'
$ws.Cells.Item(42, 5).Value = '  +5.24%  '
$ws.Cells.Item(43, 2).Value = 'WEMIXToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(43, 4).Value = '0.992'
$ws.Cells.Item(43, 5).Value = '  -0.05%  '
$ws.Cells.Item(44, 2).Value = 'MXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(44, 4).Value = '2.18'
$ws.Cells.Item(44, 5).Value = '  +0.42%  '
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).Value = '62.42'
$ws.Cells.Item(45, 5).Value = '  -0.77%  '
$ws.Cells.Item(46, 2).Value = 'RenderToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(46, 4).Value = '1.73'
$ws.Cells.Item(46, 5).Value = '  -1.53%  '
$ws.Cells.Item(47, 4).Value = '1.662.91'
$ws.Cells.Item(47, 5).Value = '  -1.72%  '
$ws.Cells.Item(48, 4).Value = '85.45'
$ws.Cells.Item(48, 5).Value = '  +0.07%  '
$ws.Cells.Item(49, 4).Value = '0.0508'
$ws.Cells.Item(49, 5).Value = '  +3.06%  '
$ws.Cells.Item(50, 5).Value = '  -1.74%  '
$ws.Cells.Item(51, 4).Value = '0.0943'
$ws.Cells.Item(51, 5).Value = '  +0.08%  '

Write-Host "applied cryptos update"
